$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$tcs = $t.ThemeColorScheme
# Office Theme colors in COM RGB order: dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink
$tcs.Item(1).RGB = 0x000000   # dk1 000000
$tcs.Item(2).RGB = 0xFFFFFF   # lt1 FFFFFF
$tcs.Item(3).RGB = 0x6A5444   # dk2 44546A -> stored reversed already handled by RGB()
